$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp caption (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 16:36"

# Countries that changed rank and therefore swapped position in the sorted table
# (their label in column A moves to a different row while row-local stats get refreshed below)
$ws.Range("A52").Value  = "Portugal"
$ws.Range("A53").Value  = "Bielorrusia"
$ws.Range("A146").Value = "Islandia"
$ws.Range("A147").Value = "Mali"
$ws.Range("A148").Value = "Botsuana"
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

# Refresh per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Muertes hoy, Muertes) for every row whose numbers changed

# Row 4
$ws.Range("B4").Value = 7781705
$ws.Range("C4").Value = 5481
$ws.Range("D4").Value = 4984897
$ws.Range("E4").Value = 2579910
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 216898

# Row 26
$ws.Range("B26").Value = 311846
$ws.Range("C26").Value = 733
$ws.Range("E26").Value = 34486
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 9660

# Row 52
$ws.Range("B52").Value = 82534
$ws.Range("C52").Value = 1278
$ws.Range("D52").Value = 51517
$ws.Range("E52").Value = 28967
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 2050

# Row 53
$ws.Range("B53").Value = 81982
$ws.Range("C53").Value = 477
$ws.Range("D53").Value = 76081
$ws.Range("E53").Value = 5021
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 880

# Row 73
$ws.Range("B73").Value = 40178
$ws.Range("C73").Value = 271
$ws.Range("D73").Value = 31710
$ws.Range("E73").Value = 7717
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 751

# Row 85
$ws.Range("B85").Value = 22445
$ws.Range("C85").Value = 1012
$ws.Range("D85").Value = 6366
$ws.Range("E85").Value = 15544
$ws.Range("G85").Value = 25
$ws.Range("H85").Value = 535

# Row 90
$ws.Range("B90").Value = 19777
$ws.Range("C90").Value = 364
$ws.Range("D90").Value = 15855
$ws.Range("E90").Value = 3147
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 775

# Row 96
$ws.Range("B96").Value = 15097
$ws.Range("C96").Value = 85
$ws.Range("E96").Value = 2959

# Row 108
$ws.Range("B108").Value = 10097
$ws.Range("C108").Value = 42
$ws.Range("D108").Value = 8914
$ws.Range("E108").Value = 1105

# Row 143
$ws.Range("B143").Value = 3617
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 2437
$ws.Range("E143").Value = 1063

# Row 146
$ws.Range("B146").Value = 3267
$ws.Range("C146").Value = 95
$ws.Range("D146").Value = 2411
$ws.Range("E146").Value = 846
$ws.Range("H146").Value = 10

# Row 147
$ws.Range("B147").Value = 3210
$ws.Range("D147").Value = 2502
$ws.Range("E147").Value = 577
$ws.Range("H147").Value = 131

# Row 148
$ws.Range("D148").Value = 834
$ws.Range("E148").Value = 2320
$ws.Range("H148").Value = 18

# Row 165
$ws.Range("B165").Value = 1360
$ws.Range("C165").Value = 5
$ws.Range("E165").Value = 33

# Row 180
$ws.Range("B180").Value = 477
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 456
$ws.Range("E180").Value = 21

# Row 215
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

# Row 216
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
